# feat: Add "Commercial" column for user assignment in full imports
#
# Adds a new column M ("Commercial") to the example import sheet so that
# rows can carry the email of the user a customer should be assigned to.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Header for the new column
$ws.Range("M1").Value = "Commercial"

# Example rows: assign the sample customers to a commercial by email
$ws.Range("M2").Value = "admin@test.com"
$ws.Range("M3").Value = "user@test.com"
